$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 3 new rows at position 60, pushing the old rows
#     60-62 down to 63-65 (their content is otherwise unchanged). ---
$ws.Rows("60:62").Insert()

# --- Step 2: update the (previously) most-recent week's rows (54-56):
#     new date 2021-10-07 (44476) and updated volume/price figures. ---
$ws.Range("D54").Value = 44476
$ws.Range("M54").Value = 45
$ws.Range("N54").Value = 30000
$ws.Range("O54").Value = 30000
$ws.Range("P54").Value = 30000
$ws.Range("S54").Value = 3000

$ws.Range("D55").Value = 44476
$ws.Range("M55").Value = 58
$ws.Range("N55").Value = 27000
$ws.Range("O55").Value = 27000
$ws.Range("P55").Value = 27000
$ws.Range("S55").Value = 2700

$ws.Range("D56").Value = 44476
$ws.Range("M56").Value = 50
$ws.Range("N56").Value = 23000
$ws.Range("O56").Value = 23000
$ws.Range("P56").Value = 23000
$ws.Range("S56").Value = 2300

# --- Step 3: rows 57-59 move their date forward one cycle
#     (44466 -> 44468); row 59 also gets new price figures. ---
$ws.Range("D57").Value = 44468

$ws.Range("D58").Value = 44468

$ws.Range("D59").Value = 44468
$ws.Range("N59").Value = 22000
$ws.Range("O59").Value = 22000
$ws.Range("P59").Value = 22000
$ws.Range("S59").Value = 2200

# --- Step 4: fill the 3 freshly inserted rows (60-62) with the
#     "week of 2021-09-27" (44466) data set. ---
$ws.Range("A60").Value = 3
$ws.Range("B60").Value = "Femacal de La Calera"
$ws.Range("C60").Value = "Coquimbo"
$ws.Range("D60").Value = 44466
$ws.Range("E60").Value = 5
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100107
$ws.Range("H60").Value = "Otros"
$ws.Range("I60").Value = 100107002
$ws.Range("J60").Value = "Chirimoya"
$ws.Range("K60").Value = "Cultivar IV Región"
$ws.Range("L60").Value = "Especial"
$ws.Range("M60").Value = 45
$ws.Range("N60").Value = 27000
$ws.Range("O60").Value = 27000
$ws.Range("P60").Value = 27000
$ws.Range("Q60").Value = "$/bandeja 10 kilos"
$ws.Range("R60").Value = "Provincia del Elquí"
$ws.Range("S60").Value = 2700
$ws.Range("T60").Value = 10

$ws.Range("A61").Value = 3
$ws.Range("B61").Value = "Femacal de La Calera"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44466
$ws.Range("E61").Value = 5
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100107
$ws.Range("H61").Value = "Otros"
$ws.Range("I61").Value = 100107002
$ws.Range("J61").Value = "Chirimoya"
$ws.Range("K61").Value = "Cultivar IV Región"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 48
$ws.Range("N61").Value = 25000
$ws.Range("O61").Value = 25000
$ws.Range("P61").Value = 25000
$ws.Range("Q61").Value = "$/bandeja 10 kilos"
$ws.Range("R61").Value = "Provincia del Elquí"
$ws.Range("S61").Value = 2500
$ws.Range("T61").Value = 10

$ws.Range("A62").Value = 3
$ws.Range("B62").Value = "Femacal de La Calera"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44466
$ws.Range("E62").Value = 5
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100107
$ws.Range("H62").Value = "Otros"
$ws.Range("I62").Value = 100107002
$ws.Range("J62").Value = "Chirimoya"
$ws.Range("K62").Value = "Cultivar IV Región"
$ws.Range("L62").Value = "Segunda"
$ws.Range("M62").Value = 40
$ws.Range("N62").Value = 23000
$ws.Range("O62").Value = 23000
$ws.Range("P62").Value = 23000
$ws.Range("Q62").Value = "$/bandeja 10 kilos"
$ws.Range("R62").Value = "Provincia del Elquí"
$ws.Range("S62").Value = 2300
$ws.Range("T62").Value = 10

# Apply the same date number-format (style index 2 in the original
# file) used by all other "Fecha" column cells.
$ws.Range("D60:D62").NumberFormat = $ws.Range("D59").NumberFormat
